$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: a new price-week is prepended (rows 70-71), pushing the
# existing history down by one pair of rows, and the oldest pair that fell
# off the bottom reappears as a freshly-dated pair at the end (rows 148-149).

# Insert two blank rows before row 70; existing rows 70-147 shift down to 72-149
$ws.Rows("70:71").Insert()

# New row 70 (Primera) for the new week
$ws.Range("A70").Value = 11
$ws.Range("B70").Value = "Vega Monumental Concepción"
$ws.Range("C70").Value = "Bíobío"
$ws.Range("D70").Value = 44705
$ws.Range("E70").Value = 8
$ws.Range("F70").Value = 100112044
$ws.Range("G70").Value = "Perejil"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 200
$ws.Range("K70").Value = 600
$ws.Range("L70").Value = 700
$ws.Range("M70").Value = 650
$ws.Range("N70").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O70").Value = "Región de Ñuble"
$ws.Range("P70").Value = 650
$ws.Range("Q70").Value = 1
$ws.Range("R70").Value = "Hortaliza"

# New row 71 (Segunda) for the new week
$ws.Range("A71").Value = 11
$ws.Range("B71").Value = "Vega Monumental Concepción"
$ws.Range("C71").Value = "Bíobío"
$ws.Range("D71").Value = 44705
$ws.Range("E71").Value = 8
$ws.Range("F71").Value = 100112044
$ws.Range("G71").Value = "Perejil"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Segunda"
$ws.Range("J71").Value = 100
$ws.Range("K71").Value = 500
$ws.Range("L71").Value = 500
$ws.Range("M71").Value = 500
$ws.Range("N71").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O71").Value = "Región de Ñuble"
$ws.Range("P71").Value = 500
$ws.Range("Q71").Value = 1
$ws.Range("R71").Value = "Hortaliza"

# New row 148 (Primera) appended at the bottom
$ws.Range("A148").Value = 11
$ws.Range("B148").Value = "Vega Monumental Concepción"
$ws.Range("C148").Value = "Bíobío"
$ws.Range("D148").Value = 44442
$ws.Range("E148").Value = 8
$ws.Range("F148").Value = 100112044
$ws.Range("G148").Value = "Perejil"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 300
$ws.Range("K148").Value = 600
$ws.Range("L148").Value = 700
$ws.Range("M148").Value = 650
$ws.Range("N148").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O148").Value = "Región de Ñuble"
$ws.Range("P148").Value = 650
$ws.Range("Q148").Value = 1
$ws.Range("R148").Value = "Hortaliza"

# New row 149 (Segunda) appended at the bottom
$ws.Range("A149").Value = 11
$ws.Range("B149").Value = "Vega Monumental Concepción"
$ws.Range("C149").Value = "Bíobío"
$ws.Range("D149").Value = 44442
$ws.Range("E149").Value = 8
$ws.Range("F149").Value = 100112044
$ws.Range("G149").Value = "Perejil"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Segunda"
$ws.Range("J149").Value = 150
$ws.Range("K149").Value = 500
$ws.Range("L149").Value = 500
$ws.Range("M149").Value = 500
$ws.Range("N149").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O149").Value = "Región de Ñuble"
$ws.Range("P149").Value = 500
$ws.Range("Q149").Value = 1
$ws.Range("R149").Value = "Hortaliza"
